# Add the new "calculator options" reference table to the "Menu Options" sheet
# (rows 10-17), documenting the new factorial operation alongside the existing ones.
# All the new cells share the same documentation font style: Calibri, size 11.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Menu Options" sheet

$newCells = @(
    @{ Cell = "A10"; Text = 'Option Number' },
    @{ Cell = "B10"; Text = 'Option Name' },
    @{ Cell = "C10"; Text = 'Function' },
    @{ Cell = "D10"; Text = 'Input Required' },
    @{ Cell = "E10"; Text = 'Output' },
    @{ Cell = "F10"; Text = 'Example' },
    @{ Cell = "G10"; Text = 'Error Handling' },
    @{ Cell = "A11"; Text = '-------------|-------------|-----------------------------------|----------------------|-----------------------|--------------------------------|--------------' },
    @{ Cell = "A12"; Text = '1' },
    @{ Cell = "B12"; Text = 'Addition' },
    @{ Cell = "C12"; Text = 'Adds two numbers together' },
    @{ Cell = "D12"; Text = 'Two floating-point...| Sum of the numbers' },
    @{ Cell = "E12"; Text = '5.5 + 3.2 = 8.7' },
    @{ Cell = "F12"; Text = 'Handles invalid input gracefully' },
    @{ Cell = "A13"; Text = '2' },
    @{ Cell = "B13"; Text = 'Subtract' },
    @{ Cell = "C13"; Text = 'Subt0rctes second number from first| Two floating-point..| Difference of the numb...| 10.0 - 4.5 = 5.5' },
    @{ Cell = "D13"; Text = 'Handles invalid input gracefully' },
    @{ Cell = "A14"; Text = '3' },
    @{ Cell = "B14"; Text = 'Multiply' },
    @{ Cell = "C14"; Text = 'Multiplies two numbers together' },
    @{ Cell = "D14"; Text = 'Two floating-point..| Product of the numbers| 3.0 * 4.0 = 12.0' },
    @{ Cell = "E14"; Text = 'Handles invalid input gracefully' },
    @{ Cell = "A15"; Text = '4' },
    @{ Cell = "B15"; Text = 'Divide' },
    @{ Cell = "C15"; Text = 'Divides first number by second' },
    @{ Cell = "D15"; Text = 'Two floating-point..| Quotient of the numb...| 15.0 / 3.0 = 5.0' },
    @{ Cell = "E15"; Text = 'Raises ValueError for division by zero; handles other errors gracefully' },
    @{ Cell = "A16"; Text = '5' },
    @{ Cell = "B16"; Text = 'Power' },
    @{ Cell = "C16"; Text = 'Raises first number to power of se...| Two floating-point..| Result of exponenti...| 2.0 ^ 3.0 = 8.0' },
    @{ Cell = "D16"; Text = 'Handles invalid input gracefully, including non-numeric inputs and negative results for even roots (raising an error or returning a complex result)' },
    @{ Cell = "A17"; Text = '6' },
    @{ Cell = "B17"; Text = 'Show History| Displays all previous calculations' },
    @{ Cell = "C17"; Text = 'None' },
    @{ Cell = "D17"; Text = 'List of calculation...' },
    @{ Cell = "E17"; Text = '"5.5 + 3.2 = 8.7", ...' },
    @{ Cell = "F17"; Text = 'Shows the history list, handles empty lists gracefully' }
)

foreach ($item in $newCells) {
    $c = $ws.Range($item.Cell)
    $c.Value = $item.Text
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
}

Write-Host "Added $($newCells.Count) cells across rows 10-17 to $($ws.Name)"